# Insert a new weekly record row before current row 63 (shifts rows 63..123 down to 64..124)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("63:63").Insert()

$ws.Range("A63").Value = 9
$ws.Range("B63").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C63").Value = "Metropolitana"
$ws.Range("D63").Value = 44778
$ws.Range("E63").Value = 13
$ws.Range("F63").Value = 100112022
$ws.Range("G63").Value = "Arveja Verde"
$ws.Range("H63").Value = "Perfection"
$ws.Range("I63").Value = "Primera"
$ws.Range("J63").Value = 29
$ws.Range("K63").Value = 40000
$ws.Range("L63").Value = 42000
$ws.Range("M63").Value = 40966
$ws.Range("N63").Value = '$/malla 25 kilos'
$ws.Range("O63").Value = "Provincia de Huasco"
$ws.Range("P63").Value = 1639
$ws.Range("Q63").Value = 25
$ws.Range("R63").Value = "Hortaliza"
